$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Hora") goes from 7 -> 8 for every data row (2-51).
# These were stored as text in the sheet, so force text with a leading
# apostrophe (Excel's text-prefix) to avoid turning them into numbers.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 7).Value = "'8"
}

# Column D ("Price") updates - numeric-looking values stored as text.
$ws.Cells.Item(2, 4).Value = "'243.89"
$ws.Cells.Item(3, 4).Value = "'23.95"
$ws.Cells.Item(4, 4).Value = "'5.168"
$ws.Cells.Item(5, 4).Value = "'0.05757"
$ws.Cells.Item(6, 4).Value = "'6.482"
$ws.Cells.Item(7, 4).Value = "'3.154"
$ws.Cells.Item(8, 4).Value = "'0.8104"
$ws.Cells.Item(9, 4).Value = "'0.8458"
$ws.Cells.Item(11, 4).Value = "'0.06954"
$ws.Cells.Item(12, 4).Value = "'0.03125"
$ws.Cells.Item(13, 4).Value = "'0.02852"
$ws.Cells.Item(14, 4).Value = "'0.09365"
$ws.Cells.Item(15, 4).Value = "'3.763"
$ws.Cells.Item(16, 4).Value = "'0.001509"
$ws.Cells.Item(17, 4).Value = "'0.04675"
$ws.Cells.Item(18, 4).Value = "'0.0005974"
$ws.Cells.Item(19, 4).Value = "'0.006140"
$ws.Cells.Item(20, 4).Value = "'0.001238"
$ws.Cells.Item(21, 4).Value = "'0.004278"
$ws.Cells.Item(22, 4).Value = "'0.00008713"
$ws.Cells.Item(23, 4).Value = "'3.502"
$ws.Cells.Item(25, 4).Value = "'0.3171"
$ws.Cells.Item(26, 4).Value = "'0.1338"
$ws.Cells.Item(27, 4).Value = "'0.1360"
$ws.Cells.Item(28, 4).Value = "'0.0002331"
$ws.Cells.Item(41, 4).Value = "'0.002975"
$ws.Cells.Item(43, 4).Value = "'0.002924"
$ws.Cells.Item(44, 4).Value = "'0.007366"
$ws.Cells.Item(45, 4).Value = "'0.00005306"
$ws.Cells.Item(47, 4).Value = "'0.3002"
$ws.Cells.Item(48, 4).Value = "'0.002277"
$ws.Cells.Item(49, 4).Value = "'0.00002101"
$ws.Cells.Item(50, 4).Value = "'0.0002001"

# Column E ("Volume(1h)") text updates.
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Write-Host "Updated symbol list"
